$wb = $excel.ActiveWorkbook

# --- 1. Update the status text "Ready for handoff" -> "In Translation" ---
# This shared string is used on the Overview sheet (columns for each locale)
# as well as on each per-locale sheet's "Status" column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the "Status" columns to fit the new, shorter text ---
# Original stored column width was ~17.216 (fit for "Ready for handoff");
# new stored column width should be ~13.410 (fit for "In Translation").
# ColumnWidth (character units) = 12.5 is the closest settable value that
# yields that narrower stored width for this engine's pixel-grid rounding.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe.Range("C1").ColumnWidth = 12.5
